$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve existing text formatting for Price/Volume columns so that
# numeric-looking strings (e.g. "0.620", "57.10") are not coerced into
# numbers and lose trailing zeros / formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.486.69"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.026.67"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.45%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.68"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.10"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -6.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.385"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0784"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.55"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.330.42"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.817"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.20"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.35"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.057.30"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.446.70"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.60"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0849"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.19"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.58"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.61"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.66%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.20"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.89"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -11.93%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.120"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0663"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.71"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.56"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.01%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.31"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.95%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.09%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0215"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.397.63"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.02"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.92"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.35"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.221.90"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.68%  "

Write-Host "Updated $($ws.Name) with latest crypto prices"
